# Update cryptocurrency price/volume data per upstream refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.794.44'
$ws.Cells.Item(2, 5).Value = '  +5.05%  '
$ws.Cells.Item(3, 4).Value = '2.463.39'
$ws.Cells.Item(3, 5).Value = '  +2.70%  '
$ws.Cells.Item(4, 5).Value = '  +0.31%  '
$ws.Cells.Item(5, 4).Value = '''159.33'
$ws.Cells.Item(5, 5).Value = '  +7.57%  '
$ws.Cells.Item(6, 4).Value = '''497.52'
$ws.Cells.Item(6, 5).Value = '  +3.50%  '
$ws.Cells.Item(7, 5).Value = '  +23.18%  '
$ws.Cells.Item(8, 5).Value = '  -0.53%  '
$ws.Cells.Item(9, 4).Value = '2.492.65'
$ws.Cells.Item(9, 5).Value = '  +4.10%  '
$ws.Cells.Item(10, 4).Value = '''6.30'
$ws.Cells.Item(10, 5).Value = '  +15.68%  '
$ws.Cells.Item(11, 5).Value = '  +5.37%  '
$ws.Cells.Item(12, 4).Value = '''0.338'
$ws.Cells.Item(12, 5).Value = '  +4.42%  '
$ws.Cells.Item(13, 5).Value = '  +1.46%  '
$ws.Cells.Item(14, 4).Value = '2.886.65'
$ws.Cells.Item(14, 5).Value = '  +2.63%  '
$ws.Cells.Item(15, 4).Value = '58.684.24'
$ws.Cells.Item(15, 5).Value = '  +4.18%  '
$ws.Cells.Item(16, 4).Value = '''21.84'
$ws.Cells.Item(16, 5).Value = '  +7.37%  '
$ws.Cells.Item(17, 5).Value = '  +2.54%  '
$ws.Cells.Item(18, 4).Value = '2.479.65'
$ws.Cells.Item(18, 5).Value = '  +3.74%  '
$ws.Cells.Item(19, 5).Value = '  +5.33%  '
$ws.Cells.Item(20, 4).Value = '''330.65'
$ws.Cells.Item(20, 5).Value = '  +4.69%  '
$ws.Cells.Item(21, 4).Value = '''10.12'
$ws.Cells.Item(21, 5).Value = '  +3.37%  '
$ws.Cells.Item(22, 4).Value = '''1.00'
$ws.Cells.Item(22, 5).Value = '  +0.06%  '
$ws.Cells.Item(23, 4).Value = '''5.98'
$ws.Cells.Item(23, 5).Value = '  +5.12%  '
$ws.Cells.Item(24, 4).Value = '''58.97'
$ws.Cells.Item(24, 5).Value = '  +3.90%  '
$ws.Cells.Item(25, 5).Value = '  +3.93%  '
$ws.Cells.Item(26, 5).Value = '  +5.33%  '
$ws.Cells.Item(27, 4).Value = '''0.993'
$ws.Cells.Item(27, 5).Value = '  -0.81%  '
$ws.Cells.Item(28, 4).Value = '2.568.14'
$ws.Cells.Item(28, 5).Value = '  +2.85%  '
$ws.Cells.Item(29, 5).Value = '  +1.57%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0809'
$ws.Cells.Item(30, 5).Value = '  +4.18%  '
$ws.Cells.Item(31, 4).Value = '''0.998'
$ws.Cells.Item(31, 5).Value = '  -0.22%  '
$ws.Cells.Item(32, 4).Value = '''19.09'
$ws.Cells.Item(32, 5).Value = '  +6.02%  '
$ws.Cells.Item(33, 4).Value = '''152.76'
$ws.Cells.Item(33, 5).Value = '  +2.85%  '
$ws.Cells.Item(34, 5).Value = '  +3.87%  '
$ws.Cells.Item(35, 5).Value = '  +9.44%  '
$ws.Cells.Item(36, 4).Value = '''3.89'
$ws.Cells.Item(36, 5).Value = '  +8.34%  '
$ws.Cells.Item(37, 5).Value = '  +6.45%  '
$ws.Cells.Item(38, 4).Value = '''0.850'
$ws.Cells.Item(38, 5).Value = '  +0.30%  '
$ws.Cells.Item(39, 4).Value = '''1.44'
$ws.Cells.Item(39, 5).Value = '  +6.83%  '
$ws.Cells.Item(40, 5).Value = '  +8.15%  '
$ws.Cells.Item(41, 4).Value = '''34.50'
$ws.Cells.Item(41, 5).Value = '  +3.17%  '
$ws.Cells.Item(42, 4).Value = '''285.82'
$ws.Cells.Item(42, 5).Value = '  +12.31%  '
$ws.Cells.Item(43, 5).Value = '  +6.87%  '
$ws.Cells.Item(44, 4).Value = '''0.610'
$ws.Cells.Item(44, 5).Value = '  +4.46%  '
$ws.Cells.Item(45, 4).Value = '''0.991'
$ws.Cells.Item(45, 5).Value = '  -0.45%  '
$ws.Cells.Item(46, 4).Value = '''0.0548'
$ws.Cells.Item(46, 5).Value = '  +1.01%  '
$ws.Cells.Item(47, 4).Value = '''0.0236'
$ws.Cells.Item(47, 5).Value = '  +5.12%  '
$ws.Cells.Item(48, 5).Value = '  +3.60%  '
$ws.Cells.Item(49, 4).Value = '''10.28'
$ws.Cells.Item(49, 5).Value = '  +0.56%  '
$ws.Cells.Item(50, 4).Value = '''0.713'
$ws.Cells.Item(50, 5).Value = '  +13.88%  '
$ws.Cells.Item(51, 4).Value = '''18.20'
$ws.Cells.Item(51, 5).Value = '  +6.66%  '
